$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters: A Rang, B Speler, C Score, D 180'ers, E 100+ finishes,
# F Totaal Score, G Aantal Darts, H 3-Darts Gemiddelde, I Totaal, J Winnaar

$rows = @(
    @{ Row=2;  A=1;  B='Burger Peach';            C=56; D=1; E=2; F=32924; G=1557; H=63.44; I=59; J=3 },
    @{ Row=3;  A=2;  B='Nick Fitzpatrick';         C=30; D=4; E=0; F=16869; G=799;  H=63.34; I=34; J=1 },
    @{ Row=4;  A=3;  B='Yannick den Daggelder';    C=24; D=0; E=0; F=15497; G=855;  H=54.38; I=24; J=1 },
    @{ Row=5;  A=4;  B='Niels van Dommelen';       C=19; D=0; E=0; F=15894; G=910;  H=52.4;  I=19; J=0 },
    @{ Row=6;  A=5;  B='Rocky Van Den Eeckhoudt';  C=13; D=0; E=1; F=11777; G=619;  H=57.08; I=14; J=0 },
    @{ Row=7;  A=6;  B='Nigel Riedel';             C=12; D=1; E=0; F=11924; G=681;  H=52.53; I=13; J=0 },
    @{ Row=8;  A=7;  B='Lukas G';                  C=11; D=1; E=0; F=9785;  G=545;  H=53.86; I=12; J=0 },
    @{ Row=9;  A=8;  B='Noah B';                   C=8;  D=2; E=0; F=8798;  G=539;  H=48.97; I=10; J=0 },
    @{ Row=10; A=9;  B='Sion Foulkes';             C=8;  D=0; E=0; F=6951;  G=454;  H=45.93; I=8;  J=0 },
    @{ Row=11; A=9;  B='joselito Vanbecelaere';    C=8;  D=0; E=0; F=4436;  G=293;  H=45.42; I=8;  J=0 },
    @{ Row=12; A=10; B='Aelbrecht Wesley';         C=5;  D=0; E=0; F=2740;  G=200;  H=41.1;  I=5;  J=0 },
    @{ Row=13; A=11; B='Constantinos Mavroudis';   C=4;  D=0; E=0; F=1991;  G=134;  H=44.57; I=4;  J=0 },
    @{ Row=14; A=13; B='Danny Littler';            C=3;  D=0; E=0; F=0;     G=0;    H=$null; I=3;  J=0 },
    @{ Row=15; A=13; B='David Wallam';             C=3;  D=0; E=0; F=1892;  G=138;  H=41.13; I=3;  J=0 },
    @{ Row=16; A=13; B='Tristan Snoep';            C=3;  D=0; E=0; F=1193;  G=96;   H=37.28; I=3;  J=0 },
    @{ Row=17; A=16; B='Arnold Van Der Vlies';     C=2;  D=0; E=0; F=2362;  G=167;  H=42.43; I=2;  J=0 },
    @{ Row=18; A=16; B='Quintin Marais';           C=2;  D=0; E=0; F=2388;  G=107;  H=66.95; I=2;  J=0 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    if ($null -eq $r.H) {
        $ws.Cells.Item($row, 8).Value = ""
    } else {
        $ws.Cells.Item($row, 8).Value = $r.H
    }
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
}
